$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.159.45"
$ws.Range("E2").Value = "  +0.73%  "
$ws.Range("D3").Value = "3.519.85"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("D5").Value = "'593.44"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").Value = "'173.85"
$ws.Range("E6").Value = "  +2.94%  "
$ws.Range("D8").Value = "'0.594"
$ws.Range("E8").Value = "  +4.08%  "
$ws.Range("E9").Value = "  +7.19%  "
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "4.129.16"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("E13").Value = "  +0.14%  "
$ws.Range("D14").Value = "'29.12"
$ws.Range("E14").Value = "  +3.65%  "
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").Value = "67.148.18"
$ws.Range("D17").Value = "3.546.48"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "'6.33"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "'14.25"
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").Value = "'395.16"
$ws.Range("E20").Value = "  +1.35%  "
$ws.Range("D21").Value = "'8.02"
$ws.Range("E21").Value = "  +0.67%  "
$ws.Range("D22").Value = "'73.11"
$ws.Range("E22").Value = "  -0.40%  "
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("E25").Value = "  -3.34%  "
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "'6.27"
$ws.Range("E29").Value = "  -2.37%  "
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("D32").Value = "'23.93"
$ws.Range("E32").Value = "  +1.53%  "
$ws.Range("D33").Value = "'7.37"
$ws.Range("E33").Value = "  -0.95%  "
$ws.Range("E34").Value = "  +1.89%  "
$ws.Range("D35").Value = "'162.98"
$ws.Range("D36").Value = "'0.899"
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("E38").Value = "  +3.71%  "
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("E40").Value = "  +4.67%  "
$ws.Range("D41").Value = "'0.0746"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "'26.42"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").Value = "'2.64"
$ws.Range("E43").Value = "  +4.25%  "
$ws.Range("D44").Value = "2.803.00"
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").Value = "'42.95"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("D47").Value = "'336.49"
$ws.Range("E47").Value = "  -5.02%  "
$ws.Range("E48").Value = "  +1.26%  "
$ws.Range("D49").Value = "'33.63"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").Value = "'6.53"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("E51").Value = "  -0.34%  "
